# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 1096
    6  = 54
    8  = 11284
    9  = 4294
    13 = 2508
    15 = 111
    17 = 168
    18 = 492
    19 = 11255
    20 = 11107
    25 = 37
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}

$wb.Save()
